# Created finished config button.
# Adds two new rows (84 and 85) to the "Translation" sheet, each describing a
# new single-use text id ("SingleUseId84" / "SingleUseId85") with the same
# default formatting metadata used by the other rows in that table
# (Default / Left / LTR / "Value from main.c: <value>").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B84").Value = "SingleUseId84"
$ws.Range("C84").Value = "Default"
$ws.Range("D84").Value = "Left"
$ws.Range("E84").Value = "LTR"
$ws.Range("F84").Value = "Value from main.c: <value>"

$ws.Range("B85").Value = "SingleUseId85"
$ws.Range("C85").Value = "Default"
$ws.Range("D85").Value = "Left"
$ws.Range("E85").Value = "LTR"
$ws.Range("F85").Value = "Value from main.c: <value>"
